$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Number"
$ws.Range("B1").Value = "Name"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Karyna"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Natasha"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Vitaliy"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Alex"

$ws.Range("B5").Select()
